# se realiza ajuste del DataExcelModels y se quitan los datos quemados del portinPrepago
#
# - D13 (portId column) changes from "732111324707277" to "732111324707278"
#   on both sheets.
# - A new row 14 is appended on both sheets, re-using the row-13 MSI/Vendedor
#   values, a new raw numeric Cedula (3043209863) and the old portId value
#   (now with a trailing space) moved down into D14.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Semilla 3", "Semilla 6")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 13: bump the "quemado" portId value.
    $ws.Range("D13").NumberFormat = "@"
    $ws.Range("D13").Value = "732111324707278"

    # New row 14: same MSI/Vendedor as row 13, a plain numeric Cedula, and
    # the previous portId value (with a trailing space) as text.
    $ws.Range("A14").NumberFormat = "@"
    $ws.Range("A14").Value = "10960370"

    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "621218573"

    $ws.Range("C14").Value = 3043209863
    $ws.Range("C14").NumberFormat = "@"

    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = "732111324707277 "
}

# Restore the active-cell selections recorded in the workbook after the edit.
$ws1 = $wb.Worksheets.Item("Semilla 3")
$ws2 = $wb.Worksheets.Item("Semilla 6")

$ws2.Activate() | Out-Null
$ws2.Range("C14:D14").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A13:B14").Select() | Out-Null
